# Ajustando os pontos no mapa utilizando
# Removes the trailing degree symbol (and stray whitespace) from the
# Latitude/Longitude text values in columns L and M, and updates the
# active cell selection to M2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to remain a text (string) cell even though the new
    # value looks numeric, matching the original authoring (t="s").
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
}

Set-TextValue "L2" "-15.793404"
Set-TextValue "M2" "-47.882317"
Set-TextValue "L4" "-15.833528"
Set-TextValue "M4" "-48.056572"
Set-TextValue "L5" "-15.68089"
Set-TextValue "M5" "-48.194262"
Set-TextValue "L6" "-15.650053"
Set-TextValue "M6" "-47.784845"
Set-TextValue "L7" "-15.618195"
Set-TextValue "M7" "-47.65557"
Set-TextValue "L8" "-15.77544"
Set-TextValue "M8" "-47.779763"
Set-TextValue "L9" "-15.871185"
Set-TextValue "M9" "-47.967994"
Set-TextValue "L10" "-15.817339"
Set-TextValue "M10" "-48.104577"
Set-TextValue "L11" "-15.823563"
Set-TextValue "M11" "-47.976816"
Set-TextValue "L12" "-15.790782"
Set-TextValue "M12" "-47.937443"
Set-TextValue "L13" "-15.876999"
Set-TextValue "M13" "-48.0881"
Set-TextValue "L14" "-16.017123"
Set-TextValue "M14" "-48.013133"
Set-TextValue "L15" "-15.903377"
Set-TextValue "M15" "-47.771774"
Set-TextValue "L16" "-15.918837"
Set-TextValue "M16" "-48.054121"
Set-TextValue "L17" "-15.839182"
Set-TextValue "M17" "-47.875534"
Set-TextValue "L18" "-15.883305"
Set-TextValue "M18" "-48.017476"
Set-TextValue "L19" "-15.734235"
Set-TextValue "M19" "-47.864158"
Set-TextValue "L20" "-15.853661"
Set-TextValue "M20" "-47.949377"
Set-TextValue "L21" "-15.841993"
Set-TextValue "M21" "-48.028121"
Set-TextValue "L22" "-15.907114"
Set-TextValue "M22" "-48.048811"
Set-TextValue "L23" "-15.800219"
Set-TextValue "M23" "-47.92439"
Set-TextValue "L24" "-15.710824"
Set-TextValue "M24" "-47.876334"
Set-TextValue "L25" "-15.902739"
Set-TextValue "M25" "-47.963322"
Set-TextValue "L26" "-15.782401"
Set-TextValue "M26" "-47.987921"
Set-TextValue "L27" "-15.650053"
Set-TextValue "M27" "-47.784845"
Set-TextValue "L28" "-15.863892"
Set-TextValue "M28" "-47.788521"
Set-TextValue "L29" "-15.748453"
Set-TextValue "M29" "-47.769694"
Set-TextValue "L30" "-15.806182"
Set-TextValue "M30" "-47.959158"
Set-TextValue "L31" "-15.813571"
Set-TextValue "M31" "-48.015604"
Set-TextValue "L32" "-15.600233"
Set-TextValue "M32" "-47.871595"
Set-TextValue "L33" "-15.812706"
Set-TextValue "M33" "-48.151735"
Set-TextValue "L34" "-15.858976"
Set-TextValue "M34" "-48.012228"

$ws.Range("M2").Select()
